$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.333.38'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.715.45'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.55'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5297'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06686'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2648'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.89'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07685'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.493'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.951.68'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.714.56'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5795'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8196'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.75'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '27.358.93'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '221.78'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.648'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.022'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.007'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.48'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("E26").Value = '  -2.32%  '
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.253'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.23'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05379'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.296'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.480'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.414'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.635'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.855'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9511'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.397'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5899'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.157.99'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01652'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.837'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.006'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8397'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.858.86'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈118'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.81'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4579'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.140'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05198'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.03%  '
